$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new cells to row 10 (Hobbies) to mirror column E into F, and mark G with justification
$ws.Range("F10").Value = "Hobbies are shown"
$ws.Range("G10").Value = "Changed source code"

# Add the two new cells to row 11 (CCAs) to mirror column E into F, and mark G with justification
$ws.Range("F11").Value = "CCAs are shown"
$ws.Range("G11").Value = "Changed source code"

# Update the active selection to N11 as in the after-state
$ws.Range("N11").Select()
